# Updates cryptos list (prices + 1h volume %) to the latest scrape.
# Generated from the authoritative cell-level diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.500.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.527.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.520.60"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.90%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.75"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.581"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.093.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "619.67"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -9.54%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.530.72"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.540.97"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.63%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.884"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.51%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.08%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.63"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.34"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.31%  "

# Row 32
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.88%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.50"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.70%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.34"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "570.23"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.80"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.32"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.139"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.37%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0444"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.406.23"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.327"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0715"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.50%  "

# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.51%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.47%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.31%  "
